# This script applies the weekly Esparragos price refresh described by the
# commit "Fruta / hortaliza, semanal": a new week of records (Fecha 45223,
# i.e. 2023-10-24) is inserted ahead of the existing history. Since the sheet
# keeps the newest week at the top of this block, every record from the old
# row 164 through the old row 184 is shifted down by three rows, the first
# three rows (164-166) get the brand-new weekly figures, and three extra rows
# (185-187) are appended at the end of the sheet to hold the records that used
# to be the last three (old rows 182-184). We simply rewrite every cell in the
# affected range (rows 164-187) with its final value rather than physically
# moving rows, which keeps styles/number formats intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 164
$ws.Cells.Item(164,1).Value = 9
$ws.Cells.Item(164,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(164,3).Value = "Metropolitana"
$ws.Cells.Item(164,4).Value = 45223
$ws.Cells.Item(164,5).Value = 13
$ws.Cells.Item(164,6).Value = 300000000
$ws.Cells.Item(164,7).Value = "Espárragos"
$ws.Cells.Item(164,8).Value = "Sin especificar"
$ws.Cells.Item(164,9).Value = "Banquete"
$ws.Cells.Item(164,10).Value = 52
$ws.Cells.Item(164,11).Value = 16000
$ws.Cells.Item(164,12).Value = 16000
$ws.Cells.Item(164,13).Value = 16000
$ws.Cells.Item(164,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(164,15).Value = "Provincia de Linares"
$ws.Cells.Item(164,16).Value = 1600
$ws.Cells.Item(164,17).Value = 10
$ws.Cells.Item(164,18).Value = "Hortaliza"

# Row 165
$ws.Cells.Item(165,1).Value = 9
$ws.Cells.Item(165,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(165,3).Value = "Metropolitana"
$ws.Cells.Item(165,4).Value = 45223
$ws.Cells.Item(165,5).Value = 13
$ws.Cells.Item(165,6).Value = 300000000
$ws.Cells.Item(165,7).Value = "Espárragos"
$ws.Cells.Item(165,8).Value = "Sin especificar"
$ws.Cells.Item(165,9).Value = "Primera"
$ws.Cells.Item(165,10).Value = 124
$ws.Cells.Item(165,11).Value = 14000
$ws.Cells.Item(165,12).Value = 14000
$ws.Cells.Item(165,13).Value = 14000
$ws.Cells.Item(165,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(165,15).Value = "Provincia de Linares"
$ws.Cells.Item(165,16).Value = 1400
$ws.Cells.Item(165,17).Value = 10
$ws.Cells.Item(165,18).Value = "Hortaliza"

# Row 166
$ws.Cells.Item(166,1).Value = 9
$ws.Cells.Item(166,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(166,3).Value = "Metropolitana"
$ws.Cells.Item(166,4).Value = 45223
$ws.Cells.Item(166,5).Value = 13
$ws.Cells.Item(166,6).Value = 300000000
$ws.Cells.Item(166,7).Value = "Espárragos"
$ws.Cells.Item(166,8).Value = "Sin especificar"
$ws.Cells.Item(166,9).Value = "Segunda"
$ws.Cells.Item(166,10).Value = 70
$ws.Cells.Item(166,11).Value = 12000
$ws.Cells.Item(166,12).Value = 12000
$ws.Cells.Item(166,13).Value = 12000
$ws.Cells.Item(166,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(166,15).Value = "Provincia de Linares"
$ws.Cells.Item(166,16).Value = 1200
$ws.Cells.Item(166,17).Value = 10
$ws.Cells.Item(166,18).Value = "Hortaliza"

# Row 167
$ws.Cells.Item(167,1).Value = 9
$ws.Cells.Item(167,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(167,3).Value = "Metropolitana"
$ws.Cells.Item(167,4).Value = 44490
$ws.Cells.Item(167,5).Value = 13
$ws.Cells.Item(167,6).Value = 300000000
$ws.Cells.Item(167,7).Value = "Espárragos"
$ws.Cells.Item(167,8).Value = "Sin especificar"
$ws.Cells.Item(167,9).Value = "Banquete"
$ws.Cells.Item(167,10).Value = 250
$ws.Cells.Item(167,11).Value = 1300
$ws.Cells.Item(167,12).Value = 1400
$ws.Cells.Item(167,13).Value = 1350
$ws.Cells.Item(167,14).Value = "$/kilo"
$ws.Cells.Item(167,15).Value = "Provincia de Linares"
$ws.Cells.Item(167,16).Value = 1350
$ws.Cells.Item(167,17).Value = 1
$ws.Cells.Item(167,18).Value = "Hortaliza"

# Row 168
$ws.Cells.Item(168,1).Value = 9
$ws.Cells.Item(168,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(168,3).Value = "Metropolitana"
$ws.Cells.Item(168,4).Value = 44490
$ws.Cells.Item(168,5).Value = 13
$ws.Cells.Item(168,6).Value = 300000000
$ws.Cells.Item(168,7).Value = "Espárragos"
$ws.Cells.Item(168,8).Value = "Sin especificar"
$ws.Cells.Item(168,9).Value = "Primera"
$ws.Cells.Item(168,10).Value = 520
$ws.Cells.Item(168,11).Value = 1100
$ws.Cells.Item(168,12).Value = 1200
$ws.Cells.Item(168,13).Value = 1150
$ws.Cells.Item(168,14).Value = "$/kilo"
$ws.Cells.Item(168,15).Value = "Provincia de Linares"
$ws.Cells.Item(168,16).Value = 1150
$ws.Cells.Item(168,17).Value = 1
$ws.Cells.Item(168,18).Value = "Hortaliza"

# Row 169
$ws.Cells.Item(169,1).Value = 9
$ws.Cells.Item(169,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(169,3).Value = "Metropolitana"
$ws.Cells.Item(169,4).Value = 44490
$ws.Cells.Item(169,5).Value = 13
$ws.Cells.Item(169,6).Value = 300000000
$ws.Cells.Item(169,7).Value = "Espárragos"
$ws.Cells.Item(169,8).Value = "Sin especificar"
$ws.Cells.Item(169,9).Value = "Segunda"
$ws.Cells.Item(169,10).Value = 160
$ws.Cells.Item(169,11).Value = 800
$ws.Cells.Item(169,12).Value = 1000
$ws.Cells.Item(169,13).Value = 900
$ws.Cells.Item(169,14).Value = "$/kilo"
$ws.Cells.Item(169,15).Value = "Provincia de Linares"
$ws.Cells.Item(169,16).Value = 900
$ws.Cells.Item(169,17).Value = 1
$ws.Cells.Item(169,18).Value = "Hortaliza"

# Row 170
$ws.Cells.Item(170,1).Value = 9
$ws.Cells.Item(170,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(170,3).Value = "Metropolitana"
$ws.Cells.Item(170,4).Value = 44859
$ws.Cells.Item(170,5).Value = 13
$ws.Cells.Item(170,6).Value = 300000000
$ws.Cells.Item(170,7).Value = "Espárragos"
$ws.Cells.Item(170,8).Value = "Sin especificar"
$ws.Cells.Item(170,9).Value = "Banquete"
$ws.Cells.Item(170,10).Value = 450
$ws.Cells.Item(170,11).Value = 1500
$ws.Cells.Item(170,12).Value = 1600
$ws.Cells.Item(170,13).Value = 1544
$ws.Cells.Item(170,14).Value = "$/kilo"
$ws.Cells.Item(170,15).Value = "Provincia de Linares"
$ws.Cells.Item(170,16).Value = 1544
$ws.Cells.Item(170,17).Value = 1
$ws.Cells.Item(170,18).Value = "Hortaliza"

# Row 171
$ws.Cells.Item(171,1).Value = 9
$ws.Cells.Item(171,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(171,3).Value = "Metropolitana"
$ws.Cells.Item(171,4).Value = 44859
$ws.Cells.Item(171,5).Value = 13
$ws.Cells.Item(171,6).Value = 300000000
$ws.Cells.Item(171,7).Value = "Espárragos"
$ws.Cells.Item(171,8).Value = "Sin especificar"
$ws.Cells.Item(171,9).Value = "Primera"
$ws.Cells.Item(171,10).Value = 700
$ws.Cells.Item(171,11).Value = 1100
$ws.Cells.Item(171,12).Value = 1200
$ws.Cells.Item(171,13).Value = 1171
$ws.Cells.Item(171,14).Value = "$/kilo"
$ws.Cells.Item(171,15).Value = "Provincia de Linares"
$ws.Cells.Item(171,16).Value = 1171
$ws.Cells.Item(171,17).Value = 1
$ws.Cells.Item(171,18).Value = "Hortaliza"

# Row 172
$ws.Cells.Item(172,1).Value = 9
$ws.Cells.Item(172,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(172,3).Value = "Metropolitana"
$ws.Cells.Item(172,4).Value = 44883
$ws.Cells.Item(172,5).Value = 13
$ws.Cells.Item(172,6).Value = 300000000
$ws.Cells.Item(172,7).Value = "Espárragos"
$ws.Cells.Item(172,8).Value = "Sin especificar"
$ws.Cells.Item(172,9).Value = "Primera"
$ws.Cells.Item(172,10).Value = 550
$ws.Cells.Item(172,11).Value = 1100
$ws.Cells.Item(172,12).Value = 1200
$ws.Cells.Item(172,13).Value = 1155
$ws.Cells.Item(172,14).Value = "$/kilo"
$ws.Cells.Item(172,15).Value = "Provincia de Linares"
$ws.Cells.Item(172,16).Value = 1155
$ws.Cells.Item(172,17).Value = 1
$ws.Cells.Item(172,18).Value = "Hortaliza"

# Row 173
$ws.Cells.Item(173,1).Value = 9
$ws.Cells.Item(173,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(173,3).Value = "Metropolitana"
$ws.Cells.Item(173,4).Value = 44883
$ws.Cells.Item(173,5).Value = 13
$ws.Cells.Item(173,6).Value = 300000000
$ws.Cells.Item(173,7).Value = "Espárragos"
$ws.Cells.Item(173,8).Value = "Sin especificar"
$ws.Cells.Item(173,9).Value = "Segunda"
$ws.Cells.Item(173,10).Value = 280
$ws.Cells.Item(173,11).Value = 900
$ws.Cells.Item(173,12).Value = 900
$ws.Cells.Item(173,13).Value = 900
$ws.Cells.Item(173,14).Value = "$/kilo"
$ws.Cells.Item(173,15).Value = "Provincia de Linares"
$ws.Cells.Item(173,16).Value = 900
$ws.Cells.Item(173,17).Value = 1
$ws.Cells.Item(173,18).Value = "Hortaliza"

# Row 174
$ws.Cells.Item(174,1).Value = 9
$ws.Cells.Item(174,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(174,3).Value = "Metropolitana"
$ws.Cells.Item(174,4).Value = 44649
$ws.Cells.Item(174,5).Value = 13
$ws.Cells.Item(174,6).Value = 300000000
$ws.Cells.Item(174,7).Value = "Espárragos"
$ws.Cells.Item(174,8).Value = "Verde"
$ws.Cells.Item(174,9).Value = "Primera"
$ws.Cells.Item(174,10).Value = 61
$ws.Cells.Item(174,11).Value = 34000
$ws.Cells.Item(174,12).Value = 34000
$ws.Cells.Item(174,13).Value = 34000
$ws.Cells.Item(174,14).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(174,15).Value = "Provincia de Linares"
$ws.Cells.Item(174,16).Value = 4250
$ws.Cells.Item(174,17).Value = 8
$ws.Cells.Item(174,18).Value = "Hortaliza"

# Row 175
$ws.Cells.Item(175,1).Value = 9
$ws.Cells.Item(175,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(175,3).Value = "Metropolitana"
$ws.Cells.Item(175,4).Value = 44649
$ws.Cells.Item(175,5).Value = 13
$ws.Cells.Item(175,6).Value = 300000000
$ws.Cells.Item(175,7).Value = "Espárragos"
$ws.Cells.Item(175,8).Value = "Verde"
$ws.Cells.Item(175,9).Value = "Segunda"
$ws.Cells.Item(175,10).Value = 43
$ws.Cells.Item(175,11).Value = 30000
$ws.Cells.Item(175,12).Value = 30000
$ws.Cells.Item(175,13).Value = 30000
$ws.Cells.Item(175,14).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(175,15).Value = "Provincia de Linares"
$ws.Cells.Item(175,16).Value = 3750
$ws.Cells.Item(175,17).Value = 8
$ws.Cells.Item(175,18).Value = "Hortaliza"

# Row 176
$ws.Cells.Item(176,1).Value = 9
$ws.Cells.Item(176,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(176,3).Value = "Metropolitana"
$ws.Cells.Item(176,4).Value = 44649
$ws.Cells.Item(176,5).Value = 13
$ws.Cells.Item(176,6).Value = 300000000
$ws.Cells.Item(176,7).Value = "Espárragos"
$ws.Cells.Item(176,8).Value = "Verde"
$ws.Cells.Item(176,9).Value = "Tercera"
$ws.Cells.Item(176,10).Value = 16
$ws.Cells.Item(176,11).Value = 27000
$ws.Cells.Item(176,12).Value = 27000
$ws.Cells.Item(176,13).Value = 27000
$ws.Cells.Item(176,14).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(176,15).Value = "Provincia de Linares"
$ws.Cells.Item(176,16).Value = 3375
$ws.Cells.Item(176,17).Value = 8
$ws.Cells.Item(176,18).Value = "Hortaliza"

# Row 177
$ws.Cells.Item(177,1).Value = 9
$ws.Cells.Item(177,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(177,3).Value = "Metropolitana"
$ws.Cells.Item(177,4).Value = 44341
$ws.Cells.Item(177,5).Value = 13
$ws.Cells.Item(177,6).Value = 300000000
$ws.Cells.Item(177,7).Value = "Espárragos"
$ws.Cells.Item(177,8).Value = "Sin especificar"
$ws.Cells.Item(177,9).Value = "Segunda"
$ws.Cells.Item(177,10).Value = 24
$ws.Cells.Item(177,11).Value = 28000
$ws.Cells.Item(177,12).Value = 30000
$ws.Cells.Item(177,13).Value = 29000
$ws.Cells.Item(177,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(177,15).Value = "Región Metropolitana"
$ws.Cells.Item(177,16).Value = 2900
$ws.Cells.Item(177,17).Value = 10
$ws.Cells.Item(177,18).Value = "Hortaliza"

# Row 178
$ws.Cells.Item(178,1).Value = 9
$ws.Cells.Item(178,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(178,3).Value = "Metropolitana"
$ws.Cells.Item(178,4).Value = 44341
$ws.Cells.Item(178,5).Value = 13
$ws.Cells.Item(178,6).Value = 300000000
$ws.Cells.Item(178,7).Value = "Espárragos"
$ws.Cells.Item(178,8).Value = "Sin especificar"
$ws.Cells.Item(178,9).Value = "Tercera"
$ws.Cells.Item(178,10).Value = 15
$ws.Cells.Item(178,11).Value = 24000
$ws.Cells.Item(178,12).Value = 26000
$ws.Cells.Item(178,13).Value = 25067
$ws.Cells.Item(178,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(178,15).Value = "Región Metropolitana"
$ws.Cells.Item(178,16).Value = 2507
$ws.Cells.Item(178,17).Value = 10
$ws.Cells.Item(178,18).Value = "Hortaliza"

# Row 179
$ws.Cells.Item(179,1).Value = 9
$ws.Cells.Item(179,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(179,3).Value = "Metropolitana"
$ws.Cells.Item(179,4).Value = 44868
$ws.Cells.Item(179,5).Value = 13
$ws.Cells.Item(179,6).Value = 300000000
$ws.Cells.Item(179,7).Value = "Espárragos"
$ws.Cells.Item(179,8).Value = "Sin especificar"
$ws.Cells.Item(179,9).Value = "Banquete"
$ws.Cells.Item(179,10).Value = 300
$ws.Cells.Item(179,11).Value = 1500
$ws.Cells.Item(179,12).Value = 1500
$ws.Cells.Item(179,13).Value = 1500
$ws.Cells.Item(179,14).Value = "$/kilo"
$ws.Cells.Item(179,15).Value = "Región Metropolitana"
$ws.Cells.Item(179,16).Value = 1500
$ws.Cells.Item(179,17).Value = 1
$ws.Cells.Item(179,18).Value = "Hortaliza"

# Row 180
$ws.Cells.Item(180,1).Value = 9
$ws.Cells.Item(180,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(180,3).Value = "Metropolitana"
$ws.Cells.Item(180,4).Value = 44868
$ws.Cells.Item(180,5).Value = 13
$ws.Cells.Item(180,6).Value = 300000000
$ws.Cells.Item(180,7).Value = "Espárragos"
$ws.Cells.Item(180,8).Value = "Sin especificar"
$ws.Cells.Item(180,9).Value = "Primera"
$ws.Cells.Item(180,10).Value = 650
$ws.Cells.Item(180,11).Value = 1200
$ws.Cells.Item(180,12).Value = 1300
$ws.Cells.Item(180,13).Value = 1262
$ws.Cells.Item(180,14).Value = "$/kilo"
$ws.Cells.Item(180,15).Value = "Provincia de Linares"
$ws.Cells.Item(180,16).Value = 1262
$ws.Cells.Item(180,17).Value = 1
$ws.Cells.Item(180,18).Value = "Hortaliza"

# Row 181
$ws.Cells.Item(181,1).Value = 9
$ws.Cells.Item(181,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(181,3).Value = "Metropolitana"
$ws.Cells.Item(181,4).Value = 44868
$ws.Cells.Item(181,5).Value = 13
$ws.Cells.Item(181,6).Value = 300000000
$ws.Cells.Item(181,7).Value = "Espárragos"
$ws.Cells.Item(181,8).Value = "Sin especificar"
$ws.Cells.Item(181,9).Value = "Primera"
$ws.Cells.Item(181,10).Value = 450
$ws.Cells.Item(181,11).Value = 1200
$ws.Cells.Item(181,12).Value = 1200
$ws.Cells.Item(181,13).Value = 1200
$ws.Cells.Item(181,14).Value = "$/kilo"
$ws.Cells.Item(181,15).Value = "Región Metropolitana"
$ws.Cells.Item(181,16).Value = 1200
$ws.Cells.Item(181,17).Value = 1
$ws.Cells.Item(181,18).Value = "Hortaliza"

# Row 182
$ws.Cells.Item(182,1).Value = 9
$ws.Cells.Item(182,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(182,3).Value = "Metropolitana"
$ws.Cells.Item(182,4).Value = 45216
$ws.Cells.Item(182,5).Value = 13
$ws.Cells.Item(182,6).Value = 300000000
$ws.Cells.Item(182,7).Value = "Espárragos"
$ws.Cells.Item(182,8).Value = "Sin especificar"
$ws.Cells.Item(182,9).Value = "Banquete"
$ws.Cells.Item(182,10).Value = 34
$ws.Cells.Item(182,11).Value = 16000
$ws.Cells.Item(182,12).Value = 16000
$ws.Cells.Item(182,13).Value = 16000
$ws.Cells.Item(182,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(182,15).Value = "Provincia de Linares"
$ws.Cells.Item(182,16).Value = 1600
$ws.Cells.Item(182,17).Value = 10
$ws.Cells.Item(182,18).Value = "Hortaliza"

# Row 183
$ws.Cells.Item(183,1).Value = 9
$ws.Cells.Item(183,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(183,3).Value = "Metropolitana"
$ws.Cells.Item(183,4).Value = 45216
$ws.Cells.Item(183,5).Value = 13
$ws.Cells.Item(183,6).Value = 300000000
$ws.Cells.Item(183,7).Value = "Espárragos"
$ws.Cells.Item(183,8).Value = "Sin especificar"
$ws.Cells.Item(183,9).Value = "Primera"
$ws.Cells.Item(183,10).Value = 70
$ws.Cells.Item(183,11).Value = 14000
$ws.Cells.Item(183,12).Value = 14000
$ws.Cells.Item(183,13).Value = 14000
$ws.Cells.Item(183,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(183,15).Value = "Provincia de Linares"
$ws.Cells.Item(183,16).Value = 1400
$ws.Cells.Item(183,17).Value = 10
$ws.Cells.Item(183,18).Value = "Hortaliza"

# Row 184
$ws.Cells.Item(184,1).Value = 9
$ws.Cells.Item(184,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(184,3).Value = "Metropolitana"
$ws.Cells.Item(184,4).Value = 45216
$ws.Cells.Item(184,5).Value = 13
$ws.Cells.Item(184,6).Value = 300000000
$ws.Cells.Item(184,7).Value = "Espárragos"
$ws.Cells.Item(184,8).Value = "Sin especificar"
$ws.Cells.Item(184,9).Value = "Segunda"
$ws.Cells.Item(184,10).Value = 52
$ws.Cells.Item(184,11).Value = 12000
$ws.Cells.Item(184,12).Value = 12000
$ws.Cells.Item(184,13).Value = 12000
$ws.Cells.Item(184,14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(184,15).Value = "Provincia de Linares"
$ws.Cells.Item(184,16).Value = 1200
$ws.Cells.Item(184,17).Value = 10
$ws.Cells.Item(184,18).Value = "Hortaliza"

# Row 185
$ws.Cells.Item(185,1).Value = 9
$ws.Cells.Item(185,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(185,3).Value = "Metropolitana"
$ws.Cells.Item(185,4).Value = 44460
$ws.Cells.Item(185,5).Value = 13
$ws.Cells.Item(185,6).Value = 300000000
$ws.Cells.Item(185,7).Value = "Espárragos"
$ws.Cells.Item(185,8).Value = "Sin especificar"
$ws.Cells.Item(185,9).Value = "Primera"
$ws.Cells.Item(185,10).Value = 250
$ws.Cells.Item(185,11).Value = 1900
$ws.Cells.Item(185,12).Value = 1900
$ws.Cells.Item(185,13).Value = 1900
$ws.Cells.Item(185,14).Value = "$/kilo"
$ws.Cells.Item(185,15).Value = "Región Metropolitana"
$ws.Cells.Item(185,16).Value = 1900
$ws.Cells.Item(185,17).Value = 1
$ws.Cells.Item(185,18).Value = "Hortaliza"

# Row 186
$ws.Cells.Item(186,1).Value = 9
$ws.Cells.Item(186,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186,3).Value = "Metropolitana"
$ws.Cells.Item(186,4).Value = 44460
$ws.Cells.Item(186,5).Value = 13
$ws.Cells.Item(186,6).Value = 300000000
$ws.Cells.Item(186,7).Value = "Espárragos"
$ws.Cells.Item(186,8).Value = "Sin especificar"
$ws.Cells.Item(186,9).Value = "Segunda"
$ws.Cells.Item(186,10).Value = 160
$ws.Cells.Item(186,11).Value = 1700
$ws.Cells.Item(186,12).Value = 1700
$ws.Cells.Item(186,13).Value = 1700
$ws.Cells.Item(186,14).Value = "$/kilo"
$ws.Cells.Item(186,15).Value = "Región Metropolitana"
$ws.Cells.Item(186,16).Value = 1700
$ws.Cells.Item(186,17).Value = 1
$ws.Cells.Item(186,18).Value = "Hortaliza"

# Row 187
$ws.Cells.Item(187,1).Value = 9
$ws.Cells.Item(187,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(187,3).Value = "Metropolitana"
$ws.Cells.Item(187,4).Value = 44460
$ws.Cells.Item(187,5).Value = 13
$ws.Cells.Item(187,6).Value = 300000000
$ws.Cells.Item(187,7).Value = "Espárragos"
$ws.Cells.Item(187,8).Value = "Sin especificar"
$ws.Cells.Item(187,9).Value = "Tercera"
$ws.Cells.Item(187,10).Value = 106
$ws.Cells.Item(187,11).Value = 1500
$ws.Cells.Item(187,12).Value = 1500
$ws.Cells.Item(187,13).Value = 1500
$ws.Cells.Item(187,14).Value = "$/kilo"
$ws.Cells.Item(187,15).Value = "Región Metropolitana"
$ws.Cells.Item(187,16).Value = 1500
$ws.Cells.Item(187,17).Value = 1
$ws.Cells.Item(187,18).Value = "Hortaliza"

# The three brand-new rows (185-187) start out with the default "General"
# number format; restore the same date/time format the Fecha column (D) uses
# everywhere else in the sheet so the stored style matches the rest of the data.
$ws.Cells.Item(185,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(186,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(187,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
